# Locate the start of the second paragraph ("Warum ist es ueberhaupt moeglich...")
# and replace everything from there through the end of the document body with
# the new, restructured content. The very first paragraph
# ("Fragen/Teilgebiete/Gliederungspunkte/Absaetze:") is left untouched in place,
# since a duplicate of it re-appears further down in the new structure.
$d = $word.ActiveDocument

$startRange = $d.Content.Duplicate
$startRange.Find.Execute("Warum ist es überhaupt möglich Haptische Schnittstellen für diese Funktion einzusetzen?") | Out-Null
$startPos = $startRange.Start

$endPos = $d.Paragraphs.Last.Range.End

$target = $d.Range($startPos, $endPos)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Aufmerksamkeit ist ein weitläufiges Feld. Deshalb ist es für die Diskussion in der Arbeit wichtig genau zu definieren, welche Arten der Aufmerksamkeit behandelt werden.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t>Fragen/Teilgebiete/Gliederungspunkte/Absätze:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Überwachungsaufgaben fordern von dem Ausführer, dass er über einen Längeren Zeitraum Informationen aufnimmt und wie der Name nahelegt überwacht. Ein einfaches Beispiel hierfür ist die Aufgabe eines Sicherheitsbeauftragten, der auf Überwachungsmonitore schaut. Auszeichnendes Merkmal ist, dass die Meiste Zeit die Meisten Informationen unverändert bleiben. Wie im Abschnitt Aufmerksamkeit gesehen, ist dieselbige als endliche Ressource </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>betrachtbar</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t xml:space="preserve">Des </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>weiteren</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> ist es evolutionär bedingt, für den Menschen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ermüdent</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> gleichbleibende Strukturen zu beobachten. Im Gegensatz dazu werden sich schnell </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>veränderende</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Informationen mit maximaler Aufmerksamkeit verfolgt. Hier bietet es sich jetzt an diese schnelle Veränderung durch haptische Aktuatoren zu simulieren.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Der menschliche Tastsinn ist evolutionär darauf ausgelegt, dass er schnell </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>ein Aufforderung</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> an das Großhirn sendet den Fokus de</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">r Augen auf eine bestimmte Körperstelle zu schieben. Wie in dem Buch </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ToDo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> beschrieben muss der Mensch in der Lage sein durch seinen Tastsinn in kürzester Zeit zu entscheiden, ob das </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>gesürte</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Gefährlich</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> oder harmlos ist. Diese Aufgabe wird zumeist so abgearbeitet, dass der Tastsinn eine Berührung </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>registiert</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. Diese bewirkt, dass sich die Aufmerksamkeit des Menschen auf den Punkt verlagert und er mehr Sinne zur Verfügung hat, um die Situation, in der er sich befindet einzuschätzen.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Somit ist das Ziel, der Kombination der Einzelnen Gebiete, dass das </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Ermüdendende</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>dauerhaufte</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Aufmerksam</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> sein des Menschen an Maschinen ausgelagert wird. Dabei registriert ein </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Sensore</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> eine Veränderung, die dann von einem entsprechenden Aufbau an Aktoren haptisch dargestellt werden.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($xml) | Out-Null

"Done"
